# login() is added with negative testcases
#
# Changes applied to Sheet1 of the Test Suite workbook:
#   1. Row 5 ("COD Order"): a new Action3 step "orderCOD" is appended in I5.
#   2. Row 6 ("Checking Filters"): Execution Flag flipped from YES to NO.
#   3. Row 7 ("Emailing Reports"): Execution Flag flipped from YES to NO.
#   4. A brand-new row 8 is appended: a negative test case "filters in test"
#      with Test Case ID "“ ”", Execution Flag "NO", Action1 "applyFilters".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) COD Order row gains an extra action in column I.
$ws.Range("I5").Value = "orderCOD"

# 2) Checking Filters row is disabled (negative test toggle).
$ws.Range("C6").Value = "NO"

# 3) Emailing Reports row is disabled (negative test toggle).
$ws.Range("C7").Value = "NO"

# 4) New negative test case row appended at the bottom of the sheet.
$ws.Range("A8").Value = "filters in test"
$ws.Range("B8").Value = "“ ”"
$ws.Range("C8").Value = "NO"
$ws.Range("D8").Value = "applyFilters"

# Matches the author's final cursor position recorded in the workbook.
$ws.Range("C5").Select() | Out-Null
